$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 963472
$ws.Range("C4").Value = 2821
$ws.Range("D4").Value = 118336
$ws.Range("E4").Value = 790779
$ws.Range("G4").Value = 101
$ws.Range("H4").Value = 54357

# Row 17 - Paises Bajos
$ws.Range("F17").Value = 934

# Row 47 - Republica Dominicana
$ws.Range("B47").Value = 6135
$ws.Range("C47").Value = 209
$ws.Range("D47").Value = 910
$ws.Range("E47").Value = 4947
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 278

# Row 117 - Kenia
$ws.Range("B117").Value = 355
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 106
$ws.Range("E117").Value = 235

# Row 166 - Guinea-Bisau
$ws.Range("B166").Value = 53
$ws.Range("C166").Value = 1
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 1
